$wb = $excel.ActiveWorkbook

# --- Work on the sheet that will become "calcAngle" (currently "Sheet3") ---
$calc = $wb.Worksheets.Item("Sheet3")

# Update the input values
$calc.Range("B2").Value = 99.896789999999996
$calc.Range("B4").Value = "xxxx"
$calc.Range("B5").Value = 33

# Highlight fills - green fill first (becomes fillId 2), then light blue (becomes fillId 3),
# matching the order they were introduced upstream.
$calc.Range("D15").Interior.Color = 5296274
$calc.Range("B1").Interior.Color = 15261367
$calc.Range("B2").Interior.Color = 15261367
$calc.Range("B3").Interior.Color = 15261367
$calc.Range("B5").Interior.Color = 15261367

# Rename + move to the first tab position
$calc.Name = "calcAngle"
$calc.Move($wb.Worksheets.Item(1))

# Make it the active/selected tab with B3 selected
$calcAngle = $wb.Worksheets.Item("calcAngle")
$calcAngle.Select()
$calcAngle.Range("B3").Select()

# --- Update selections on the other two sheets ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Select()
$sheet1.Range("B5").Select()

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Select()
$sheet2.Range("B4").Select()

# Leave calcAngle as the active sheet/tab
$calcAngle = $wb.Worksheets.Item("calcAngle")
$calcAngle.Select()
